$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Qminus1)
$ws.Range("B2").Value = 0.01965582767097413
$ws.Range("C2").Value = 0.9028368076043077
$ws.Range("D2").Value = 1.487274831470099
$ws.Range("E2").Value = 1.219538778173986
$ws.Range("F2").Value = 1.23381140640422
$ws.Range("G2").Value = 43

# Update existing row 3 (Q0)
$ws.Range("B3").Value = 0.09922626907716836
$ws.Range("C3").Value = 1.334569147596012
$ws.Range("D3").Value = 4.031143965198564
$ws.Range("E3").Value = 2.007770894598924
$ws.Range("F3").Value = 2.013015443340949
$ws.Range("G3").Value = 131

# Add new row 4 (Q1)
$ws.Range("A4").Value = "Q1"
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4").Value = 0.1165997114055807
$ws.Range("C4").Value = 1.458584920587325
$ws.Range("D4").Value = 8.577350374430635
$ws.Range("E4").Value = 2.928711384624752
$ws.Range("F4").Value = 2.950278682772432
$ws.Range("G4").Value = 62
